# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45189 to 45190 (i.e. bump the date by one day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 224
}

$startCell = $ws.Cells.Item(2, 3)
$endCell = $ws.Cells.Item($lastRow, 3)
$range = $ws.Range($startCell, $endCell)

foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45189) {
        $cell.Value = 45190
    }
}
